# Act graficos y tablas web pob
# - Rename sheets: "Datos" -> "Data", "Ficha técnica" -> "Metadata"
# - "Data" sheet: reorder the yearly rows so the table runs from 2020 down
#   to 2000 (descending) instead of 2000 up to 2020 (ascending), keeping the
#   (year, value) pairs intact. Implemented as a set of row swaps using
#   Range.Copy so the year cells keep their original (text) cell type
#   instead of being re-typed as numbers via a plain Value assignment.
# - "Metadata" sheet: relabel the field-name column to lower-case technical
#   keys, reorder a couple of rows, add an "observaciones" row and a final
#   attribution row.

$wb = $excel.ActiveWorkbook

$dataSheet = $wb.Worksheets.Item("Datos")
$metaSheet = $wb.Worksheets.Item("Ficha técnica")

# ---------------------------------------------------------------------
# 1) Rename the sheets
# ---------------------------------------------------------------------
$dataSheet.Name = "Data"
$metaSheet.Name = "Metadata"

# ---------------------------------------------------------------------
# 2) "Data" sheet - reverse the 21 data rows (rows 2..22) so that the
#    table goes from 2020 (top) to 2000 (bottom). Swap row r with row
#    (24 - r) for r = 2..11 (row 12 is the middle row and stays put).
#    A scratch row far below the data (row 1000) is used as a holding
#    area so Range.Copy can be used for every move, preserving the
#    original cell type (text for column A, number for column B).
# ---------------------------------------------------------------------
$scratchRow = 1000
$firstRow = 2
$lastRow = 22

for ($r = $firstRow; $r -lt ($firstRow + $lastRow) / 2; $r++) {
    $mirrorRow = $firstRow + $lastRow - $r

    $srcRange = $dataSheet.Range("A" + $r + ":B" + $r)
    $mirrorRange = $dataSheet.Range("A" + $mirrorRow + ":B" + $mirrorRow)
    $scratchRange = $dataSheet.Range("A" + $scratchRow + ":B" + $scratchRow)

    $srcRange.Copy($scratchRange)
    $mirrorRange.Copy($srcRange)
    $scratchRange.Copy($mirrorRange)
    $scratchRange.Clear()
}

# ---------------------------------------------------------------------
# 3) "Metadata" sheet - relabel / reorder rows and append new ones.
#    Row 1 (A1 empty / B1 " ") is untouched - it doesn't change.
# ---------------------------------------------------------------------
$metaSheet.Cells.Item(2, 1).Value = "nomindicador"
$metaSheet.Cells.Item(2, 2).Value = "Razón de mortalidad materna por 100000 nacidos vivos"

$metaSheet.Cells.Item(3, 1).Value = "derecho"
$metaSheet.Cells.Item(3, 2).Value = "Salud"

$metaSheet.Cells.Item(4, 1).Value = "conindicador"
$metaSheet.Cells.Item(4, 2).Value = "Mortalidad materna"

$metaSheet.Cells.Item(5, 1).Value = "tipoind"
$metaSheet.Cells.Item(5, 2).Value = "Resultados"

$metaSheet.Cells.Item(6, 1).Value = "definicion"
$metaSheet.Cells.Item(6, 2).Value = "El indicador mide la cantidad anual de defunciones maternas cada 100000 nacidos vivos."

$metaSheet.Cells.Item(7, 1).Value = "calculo"
$metaSheet.Cells.Item(7, 2).Value = "Para cada año calcular: (Cantidad de defunciones maternas ocurridas en el año acaecido / Cantidad de nacidos vivos en el año acaecido)*100000"

$metaSheet.Cells.Item(8, 1).Value = "observaciones"
$metaSheet.Cells.Item(8, 2).Value = "Sin observaciones"

$metaSheet.Cells.Item(9, 1).Value = "cita"
$metaSheet.Cells.Item(9, 2).Value = "UMAD con base en Estadísticas Vitales - MSP"

$metaSheet.Cells.Item(10, 1).Value = "Mirador DESCA - UMAD/FCS – INDDHH"
$metaSheet.Cells.Item(10, 2).Value = " "
